# "Generate Report for Handback" - record that the handback (target ->
# translated xliff) files were produced / received for both locales, and
# flag the two files as "Handed back: in sync with en-US" everywhere their
# status is shown.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$file1Name = "2ae8ff8a-1748-42c9-acdf-e0c683c8ea29.md"
$file2Name = "d4f84198-bb66-4391-8562-8ebfbb5dfee8.md"

$file1Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7a88f8919794e812e47a06213efdb2b1bc4ae27/e2e/2ae8ff8a-1748-42c9-acdf-e0c683c8ea29.md"
$file2Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c7a88f8919794e812e47a06213efdb2b1bc4ae27/e2e/d4f84198-bb66-4391-8562-8ebfbb5dfee8.md"

# ---------------------------------------------------------------------
# Overview sheet: update the per-locale status cells and widen the two
# locale status columns so the longer text fits.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Range("E1").ColumnWidth = 29.166666666666668
$wsOverview.Range("F1").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Per locale detail sheets (zh-cn, de-de): fill in the "Latest Target
# File" / "Latest Handback File" / "Latest Handback DateTime" columns now
# that the handback has happened, update the Status column, add
# hyperlinks for the newly filled in target-file column, and widen the
# columns that now hold full file names.
# ---------------------------------------------------------------------
$zhHandbackDateTime = "2016-10-13 14:32:45"
$deHandbackDateTime = "2016-10-13 14:33:03"

$locales = @(
    @{ Name = "zh-cn"; HandbackDateTime = $zhHandbackDateTime; Suffix = "zh-cn" },
    @{ Name = "de-de"; HandbackDateTime = $deHandbackDateTime; Suffix = "de-de" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    $file1Xlf = "2ae8ff8a-1748-42c9-acdf-e0c683c8ea29.1456c31ceb946685b5d5dd0867d7ea9f2d6ad205." + $locale.Suffix + ".xlf"
    $file2Xlf = "d4f84198-bb66-4391-8562-8ebfbb5dfee8.ea4e64bf4a1250c60c5a92e51c8ad1ef4075776e." + $locale.Suffix + ".xlf"

    # Status column
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Target File (column I) and Latest Handback File (column J)
    $ws.Range("I2").Value = $file1Name
    $ws.Range("J2").Value = $file1Xlf
    $ws.Range("I3").Value = $file2Name
    $ws.Range("J3").Value = $file2Xlf

    # Latest Handback DateTime (column K)
    $ws.Range("K2").Value = $locale.HandbackDateTime
    $ws.Range("K3").Value = $locale.HandbackDateTime

    # Rebuild the hyperlinks collection so link order/ids follow the row
    # order: A2, I2, A3, I3.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $file1Url, [Type]::Missing, [Type]::Missing, $file1Name)
    $ws.Hyperlinks.Add($ws.Range("I2"), $file1Url, [Type]::Missing, [Type]::Missing, $file1Name)
    $ws.Hyperlinks.Add($ws.Range("A3"), $file2Url, [Type]::Missing, [Type]::Missing, $file2Name)
    $ws.Hyperlinks.Add($ws.Range("I3"), $file2Url, [Type]::Missing, [Type]::Missing, $file2Name)

    # Restore/apply the same hyperlink look (single underline, blue font)
    # on all four anchor cells.
    $hyperlinkColor = [System.Drawing.ColorTranslator]::ToOle([System.Drawing.Color]::FromArgb(0x64, 0x95, 0xED))
    foreach ($cellRef in @("A2", "I2", "A3", "I3")) {
        $ws.Range($cellRef).Font.Underline = [Microsoft.Office.Interop.Excel.XlUnderlineStyle]::xlUnderlineStyleSingle
        $ws.Range($cellRef).Font.Color = $hyperlinkColor
    }

    # Widen Status (C) and the two file-name columns (I, J) so the longer
    # values fit.
    $ws.Range("C1").ColumnWidth = 29.166666666666668
    $ws.Range("I1").ColumnWidth = 39.166666666666664
    $ws.Range("J1").ColumnWidth = 39.166666666666664
}
